# Attendance report sync: reorders "Recorded By" names, refreshes the
# class-level missing/pending session counters, and flips the twelve
# still-outstanding "SURGERY SEMINAR/SLIDE" session-7 rows (one per B1
# group) from "Pending" to "Not Recorded" now that their due date has
# passed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Recorded By" cells: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# ---------------------------------------------------------------------
$gRows = @(2,3,4,5,6,7,8,16,17,22,23,24,26,29,37,38,43,44,45,47,50,58,59,64,65,66,68,71,79,80,85,86,87,88,89,90,91,99,100,105,106,107,108,109,110,111,119,120,125,126,127,128,129,130,131,139,140,145,146,147,148,149,150,151,159,160,165,166,167,168,169,170,171,179,180,185,186,187,189,192,200,201,206,207,208,210,213,221,222,227,228,229,231,234,242,243)

foreach ($r in $gRows) {
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------
# 2) Class Statistics: Missing Sessions / Pending Sessions totals
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 66
$ws.Range("L8").Value = 60

# ---------------------------------------------------------------------
# 3) Group Statistics (rows 15-26, one per B1 group): Missing +1, Pending -> 5
# ---------------------------------------------------------------------
$groupStatsRows = @{
    15 = 4
    16 = 5
    17 = 5
    18 = 5
    19 = 4
    20 = 4
    21 = 4
    22 = 4
    23 = 4
    24 = 5
    25 = 5
    26 = 5
}
foreach ($r in $groupStatsRows.Keys) {
    $ws.Range("P$r").Value = $groupStatsRows[$r]
    $ws.Range("Q$r").Value = 5
}

# ---------------------------------------------------------------------
# 4) Session-7 "SURGERY SEMINAR/SLIDE" rows per group: now overdue, so
#    their status flips from "Pending" (yellow, style index 6) to
#    "Not Recorded" (pink, style index 5) - reuse the formatting from an
#    existing "Not Recorded" row (row 18) instead of inventing a new style.
# ---------------------------------------------------------------------
$pendingToNotRecordedRows = @(20,41,62,83,103,123,143,163,183,204,225,246)

$formatSource = $ws.Range("A18:I18")
$formatSource.Copy()
foreach ($r in $pendingToNotRecordedRows) {
    $dst = $ws.Range("A" + $r + ":I" + $r)
    $dst.PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}
$excel.CutCopyMode = 0
